$d = $word.ActiveDocument

$replacements = @(
    @("492÷9=", "925÷7="),
    @("493÷7=", "658÷9="),
    @("791÷9=", "706÷9="),
    @("815÷8=", "360÷9="),
    @("230÷6=", "144÷2="),
    @("463÷9=", "880÷9="),
    @("729÷8=", "344÷9="),
    @("966÷4=", "478÷6="),
    @("120÷4=", "826÷8="),
    @("383÷9=", "975÷9="),
    @("582÷9=", "662÷6="),
    @("494÷9=", "276÷2="),
    @("266÷6=", "112÷3="),
    @("718÷9=", "682÷4="),
    @("759÷3=", "400÷2="),
    @("562÷3=", "217÷3="),
    @("977÷7=", "720÷2="),
    @("238÷2=", "471÷5="),
    @("769÷4=", "413÷8="),
    @("865÷3=", "785÷4="),
    @("275÷4=", "741÷5="),
    @("650÷2=", "119÷5="),
    @("572÷6=", "641÷9="),
    @("362÷8=", "400÷2="),
    @("173÷7=", "994÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
